$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "Panel"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Panel"

# Insert a new row at the top of the sheet, shifting existing data down
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Panel"
$ws.Range("B1").Value = "Genes"

# Update the selection/active cell to B1 (mirrors the view state in the diff)
$ws.Range("B1").Select()
